$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 255, shifting existing rows 255-264 down to 256-265,
# and copying formatting (incl. the date-style column D) from the row above.
$ws.Rows.Item(255).Insert()

# Populate the newly inserted row 255 with the new data record.
$ws.Cells.Item(255, 1).Value = 3
$ws.Cells.Item(255, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(255, 3).Value = "Coquimbo"
$ws.Cells.Item(255, 4).Value = Get-Date -Year 2021 -Month 11 -Day 9 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(255, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(255, 5).Value = 5
$ws.Cells.Item(255, 6).Value = 100112017
$ws.Cells.Item(255, 7).Value = "Apio"
$ws.Cells.Item(255, 8).Value = "Americana (o)"
$ws.Cells.Item(255, 9).Value = "Primera"
$ws.Cells.Item(255, 10).Value = 130
$ws.Cells.Item(255, 11).Value = 9000
$ws.Cells.Item(255, 12).Value = 9000
$ws.Cells.Item(255, 13).Value = 9000
$ws.Cells.Item(255, 14).Value = "`$/docena de matas"
$ws.Cells.Item(255, 15).Value = "Pan de Azúcar"
$ws.Cells.Item(255, 16).Value = 1500
$ws.Cells.Item(255, 17).Value = 6
$ws.Cells.Item(255, 18).Value = "Hortaliza"
